$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Version & History" ---
# Add a new V1.4 change-history entry below the existing V1.3 row (row 8),
# re-using the same row formatting (copy row 8 -> row 9) and then
# overwriting the cells that actually differ.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A8:E8").Copy($ws1.Range("A9:E9"))
$ws1.Range("A9").Value = "V1.4"
$ws1.Range("B9").Value = "Realized that I've put infomations about some future features, so temporarily removed them from the communications matrix."
# C9 (author), D9 (date) and E9 (status) stay the same as row 8, which the
# copy above already took care of.

# --- Sheet 2: "CommunicationMatrix" ---
# Remove the accidentally added future-feature signal rows (ACC, TSR, PP,
# LKS -> rows 16-19), clearing their contents and formatting. Only the
# "Additional information" helper cells in columns G and J keep their
# (now-empty) formatting, matching the rest of the table's empty rows.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A16:F19").Clear()
$ws2.Range("H16:I19").Clear()
$ws2.Range("K16:K19").Clear()
$ws2.Range("G16:G19").ClearContents()
$ws2.Range("J16:J19").ClearContents()
$ws2.Rows("16:19").AutoFit()
